$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.292.86"
$ws.Range("D3").Value = "3.495.35"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "589.18"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "133.77"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "0.486"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +6.30%  "
$ws.Range("E10").Value = "  +0.02%  "
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "0.388"
$ws.Range("E11").Value = "  +3.04%  "
$ws.Range("D12").Value = "4.090.74"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "3.495.70"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "64.250.81"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "25.31"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = "10.05"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "5.78"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "13.54"
$ws.Range("E20").Value = "  -0.71%  "
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "386.55"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "0.581"
$ws.Range("E22").Value = "  +2.66%  "
$ws.Range("D23").Value = "3.634.71"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "74.15"
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -0.52%  "
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "0.0000115"
$ws.Range("E27").Value = "  +2.13%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "7.37"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("E31").Value = "  -0.91%  "
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = "8.16"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("E33").Value = "  +3.65%  "
$ws.Range("D34").Value = "3.524.20"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "23.30"
$ws.Range("E36").Value = "  -0.72%  "
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = "5.33"
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "165.61"
$ws.Range("E40").Value = "  +2.15%  "
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "0.0786"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +0.50%  "
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "24.44"
$ws.Range("E45").Value = "  -4.58%  "
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("D48").Value = "2.428.53"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = "6.83"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("E50").Value = "  +1.48%  "
$ws.Range("E51").Value = "  -0.42%  "

$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(49,4).Style = "Normal"
